# Update calc_reactions test-case inputs for span 2 (ltr, x = 125).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Pt (load)
$ws.Range("C2").Value = 776
# xt (load position)
$ws.Range("C3").Value = 62.175258
# span_begin
$ws.Range("C4").Value = 0
# span_end
$ws.Range("C5").Value = 100

# Move the active selection from J4 to C7 (as recorded in the saved view state).
$ws.Range("C7").Select()

# Restore the workbook tab-bar ratio (was 600/1000, now 555/1000).
$excel.ActiveWindow.TabRatio = 0.555
